$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("B1").Value = "annotation"

# Update data rows 2-3 with new values
$ws.Range("A2").Value = "Neu.1.5.47"
$ws.Range("B2").Value = "Neu.CSFcN.0"
$ws.Range("A3").Value = "Neu.5.81"
$ws.Range("B3").Value = "Neu.Epend.0"

# Remove old rows 4-13 which are no longer part of the data
$ws.Range("A4:B13").ClearContents()
